$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need an explicit Text
# number format first, otherwise Excel COM auto-converts the assigned
# string into a numeric value (losing trailing zeros / changing type).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range('D2').Value = '57.188.01'
$ws.Range('E2').Value = '  -0.11%  '
$ws.Range('D3').Value = '2.415.31'
$ws.Range('E3').Value = '  -3.89%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '489.17'
$ws.Range('E5').Value = '  -1.46%  '
$ws.Range('D6').Value = '153.66'
$ws.Range('E6').Value = '  -0.22%  '
$ws.Range('B7').Value = 'USDC'
$ws.Range('C7').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D7').Value = '0.996'
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('B8').Value = 'XRP'
$ws.Range('C8').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D8').Value = '0.614'
$ws.Range('E8').Value = '  +18.58%  '
$ws.Range('D9').Value = '2.433.44'
$ws.Range('E9').Value = '  -3.79%  '
$ws.Range('D10').Value = '6.10'
$ws.Range('E10').Value = '  +5.80%  '
$ws.Range('D11').Value = '0.0999'
$ws.Range('E11').Value = '  -1.57%  '
$ws.Range('D12').Value = '0.333'
$ws.Range('E12').Value = '  -1.94%  '
$ws.Range('E13').Value = '  +1.23%  '
$ws.Range('D14').Value = '2.838.29'
$ws.Range('E14').Value = '  -3.79%  '
$ws.Range('D15').Value = '57.195.26'
$ws.Range('E15').Value = '  -0.27%  '
$ws.Range('D16').Value = '20.58'
$ws.Range('E16').Value = '  -3.97%  '
$ws.Range('E17').Value = '  -4.08%  '
$ws.Range('D18').Value = '2.434.91'
$ws.Range('E18').Value = '  -3.78%  '
$ws.Range('E19').Value = '  +1.34%  '
$ws.Range('D20').Value = '323.99'
$ws.Range('E20').Value = '  -0.24%  '
$ws.Range('D21').Value = '10.00'
$ws.Range('E21').Value = '  -3.41%  '
$ws.Range('D22').Value = '0.999'
$ws.Range('E22').Value = '  -0.14%  '
$ws.Range('D23').Value = '5.93'
$ws.Range('E23').Value = '  -0.40%  '
$ws.Range('D24').Value = '57.85'
$ws.Range('E24').Value = '  -1.21%  '
$ws.Range('E25').Value = '  -1.89%  '
$ws.Range('D26').Value = '0.994'
$ws.Range('E26').Value = '  -0.17%  '
$ws.Range('E27').Value = '  -2.88%  '
$ws.Range('D28').Value = '2.522.39'
$ws.Range('E28').Value = '  -3.68%  '
$ws.Range('D29').Value = '7.28'
$ws.Range('E29').Value = '  -4.85%  '
$ws.Range('D30').Value = '0.0₃0784'
$ws.Range('E30').Value = '  -5.66%  '
$ws.Range('E31').Value = '  +0.00%  '
$ws.Range('D32').Value = '150.88'
$ws.Range('E32').Value = '  -0.53%  '
$ws.Range('D33').Value = '18.65'
$ws.Range('E33').Value = '  +1.14%  '
$ws.Range('E34').Value = '  -1.25%  '
$ws.Range('D35').Value = '5.30'
$ws.Range('E35').Value = '  +0.06%  '
$ws.Range('D36').Value = '1.16'
$ws.Range('E36').Value = '  -0.77%  '
$ws.Range('D37').Value = '3.76'
$ws.Range('E37').Value = '  -2.35%  '
$ws.Range('D38').Value = '0.820'
$ws.Range('E38').Value = '  -9.71%  '
$ws.Range('E39').Value = '  +6.76%  '
$ws.Range('D40').Value = '34.01'
$ws.Range('E40').Value = '  -1.24%  '
$ws.Range('B41').Value = 'Bittensor'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D41').Value = '282.59'
$ws.Range('E41').Value = '  +4.94%  '
$ws.Range('B42').Value = 'Filecoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D42').Value = '3.52'
$ws.Range('E42').Value = '  -0.79%  '
$ws.Range('B43').Value = 'Stacks'
$ws.Range('C43').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D43').Value = '1.37'
$ws.Range('E43').Value = '  -3.18%  '
$ws.Range('D44').Value = '0.994'
$ws.Range('E44').Value = '  +0.02%  '
$ws.Range('D45').Value = '0.598'
$ws.Range('E45').Value = '  -3.82%  '
$ws.Range('D46').Value = '0.0530'
$ws.Range('E46').Value = '  -6.02%  '
$ws.Range('D47').Value = '10.21'
$ws.Range('E47').Value = '  +0.06%  '
$ws.Range('D48').Value = '0.0228'
$ws.Range('E48').Value = '  -1.52%  '
$ws.Range('D49').Value = '4.59'
$ws.Range('E49').Value = '  -6.88%  '
$ws.Range('D50').Value = '1.898.27'
$ws.Range('E50').Value = '  -0.20%  '
$ws.Range('D51').Value = '17.58'
$ws.Range('E51').Value = '  -2.99%  '
